$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Update the "compatible with" blurb: split it into two italic runs
#    ("This sample is compatible with the " / "Windows 10 Fall Creators
#    Update SDK (16299)"), merge the following (bookmark-only) paragraph
#    into it, and leave a new blank paragraph behind where that
#    paragraph used to end.
# ---------------------------------------------------------------------

$p2 = $d.Paragraphs.Item(2)
$p3 = $d.Paragraphs.Item(3)
$mergeRange = $d.Range($p2.Range.Start, $p3.Range.End)

$mergedXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body>' + `
  '<w:p w:rsidR="00A06E4B" w:rsidRDefault="00A06E4B" w:rsidP="00A06E4B">' + `
  '<w:pPr><w:rPr><w:i/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">This sample is compatible with the </w:t></w:r>' + `
  '<w:r><w:rPr><w:i/></w:rPr><w:t>Windows 10 Fall Creators Update SDK (16299)</w:t></w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
  '</w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$mergeRange.InsertXML($mergedXml)

# Leave a fresh, empty paragraph right after the merged one.
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertParagraphAfter()

# The freshly-inserted paragraph inherits the italic run properties of
# its predecessor; reset it back down to a genuinely blank paragraph.
$p3 = $d.Paragraphs.Item(3)
$blankXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData>' + `
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p3.Range.InsertXML($blankXml)

# ---------------------------------------------------------------------
# 2) Turn on even/default header & footer content (in addition to the
#    existing first-page header/footer). Clearing the text of the
#    primary (default) header is enough for Word to mint the even +
#    default header/footer parts and renumber the old table-based
#    first-page header/footer out of the way.
# ---------------------------------------------------------------------

$sec = $d.Sections(1)
$primaryHeader = $sec.Headers(1)
$primaryHeader.Range.Text = ""
